# Atualização de bases das ligas, do dia: 02-03-2024 às 08:34
# This script:
#  1) Swaps the content of rows 15/16 (match data got re-ordered upstream)
#  2) Swaps the content of rows 85/86 (match data got re-ordered upstream)
#  3) Appends 7 new upcoming fixtures as rows 151-157

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) & 2): swap full row contents (columns B, F..AC) between paired rows,
#    leaving the running index in column A (and the identical C/D/E values)
#    untouched.
# ---------------------------------------------------------------------------
function Swap-Rows($rowA, $rowB) {
    $rangeA = "B" + $rowA + ":AC" + $rowA
    $rangeB = "B" + $rowB + ":AC" + $rowB
    $valsA = $ws.Range($rangeA).Value()
    $valsB = $ws.Range($rangeB).Value()
    $ws.Range($rangeA).Value = $valsB
    $ws.Range($rangeB).Value = $valsA
}

Swap-Rows 15 16
Swap-Rows 85 86

# ---------------------------------------------------------------------------
# 3) Append the 7 new fixture rows (151-157). These are future matches so
#    only pre-match data (id, teams, opening/closing odds) is known; the
#    score (H/I), result (J) and post-match odds (AB/AC) are left blank.
# ---------------------------------------------------------------------------

# Copy the formatting (number formats, fonts, borders, alignment) of the
# last existing data row down onto the new rows first. Only the columns
# that actually carry a non-default style (A: bold/bordered id, E: date
# number format) are touched so that we don't create empty placeholder
# cells in columns that have no data on these (still unplayed) fixtures.
$ws.Range("A150").Copy()
$ws.Range("A151:A157").PasteSpecial(-4122)
$ws.Range("E150").Copy()
$ws.Range("E151:E157").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$newRows = @(
    @{ Row=151; A=149; B=6992674; F="Buriram United";      G="Chiangrai Utd";     E=45353.35416666666;  K=2.15;  L=3.4; M=2.875; N=2;     O=3.4;  P=3.3;  Q=-0.5;  R=2.025; S=1.775; T=2.5;  U=1.825; V=1.975 },
    @{ Row=152; A=150; B=6992680; F="Lamphun Warrior FC";  G="BG Pathum United";   E=45353.375;          K=2.7;   L=3.4; M=2.25;  N=2.8;   O=3.4;  P=2.2;  Q=0.25;  R=1.775; S=2.025; T=2.75; U=1.9;   V=1.9   },
    @{ Row=153; A=151; B=6995901; F="Police Tero FC";      G="Nakhon Pathom FC";   E=45353.41666666666;  K=2.15;  L=3.3; M=2.9;   N=2.25;  O=3.25; P=2.8;  Q=-0.25; R=2.025; S=1.775; T=2.75; U=1.975; V=1.825 },
    @{ Row=154; A=152; B=6992678; F="Muang Thong United";  G="Prachuap FC";        E=45354.3125;         K=1.95;  L=3.6; M=3.2;   N=1.85;  O=3.75; P=3.4;  Q=-0.5;  R=1.875; S=1.925; T=2.75; U=1.8;   V=2      },
    @{ Row=155; A=153; B=6992675; F="Sukhothai FC";        G="Bangkok United";     E=45354.33333333334;  K=1.833; L=3.6; M=3.5;   N=1.727; O=3.6;  P=4;    Q=-0.75; R=1.975; S=1.825; T=2.5;  U=1.825; V=1.975 },
    @{ Row=156; A=154; B=6992679; F="Uthai Thani FC";      G="Chonburi";           E=45354.375;          K=1.95;  L=3.5; M=3.25;  N=2.1;   O=3.4;  P=3;    Q=-0.25; R=1.9;   S=1.9;   T=2.75; U=1.975; V=1.825 },
    @{ Row=157; A=155; B=6992677; F="Trat FC";             G="Ratchaburi FC";      E=45355.35416666666;  K=2.875; L=3.4; M=2.15;  N=3.6;   O=3.5;  P=1.85; Q=0.5;   R=1.85;  S=1.95;  T=2.5;  U=1.85;  V=1.95  }
)

foreach ($r in $newRows) {
    $rn = $r.Row
    $ws.Cells.Item($rn, 1).Value  = $r.A                          # A - running index
    $ws.Cells.Item($rn, 2).Value  = $r.B                          # B - match id
    $ws.Cells.Item($rn, 3).Value  = "Thailand Premier League"     # C - Div
    $ws.Cells.Item($rn, 4).Value  = "Thailand Premier League"     # D - Div Original Name
    $ws.Cells.Item($rn, 5).Value  = $r.E                          # E - Date
    $ws.Cells.Item($rn, 6).Value  = $r.F                          # F - HomeTeam
    $ws.Cells.Item($rn, 7).Value  = $r.G                          # G - AwayTeam
    $ws.Cells.Item($rn, 11).Value = $r.K                          # K - oddH_op
    $ws.Cells.Item($rn, 12).Value = $r.L                          # L - oddD_op
    $ws.Cells.Item($rn, 13).Value = $r.M                          # M - oddA_op
    $ws.Cells.Item($rn, 14).Value = $r.N                          # N - oddH
    $ws.Cells.Item($rn, 15).Value = $r.O                          # O - oddD
    $ws.Cells.Item($rn, 16).Value = $r.P                          # P - oddA
    $ws.Cells.Item($rn, 17).Value = $r.Q                          # Q - Ah
    $ws.Cells.Item($rn, 18).Value = $r.R                          # R - oddAHH
    $ws.Cells.Item($rn, 19).Value = $r.S                          # S - oddAHA
    $ws.Cells.Item($rn, 20).Value = $r.T                          # T - AhOU
    $ws.Cells.Item($rn, 21).Value = $r.U                          # U - oddAHOver
    $ws.Cells.Item($rn, 22).Value = $r.V                          # V - oddAHUnder
    $ws.Cells.Item($rn, 23).Value = 0                             # W - PLH
    $ws.Cells.Item($rn, 24).Value = 0                             # X - PLD
    $ws.Cells.Item($rn, 25).Value = 0                             # Y - PLA
    $ws.Cells.Item($rn, 26).Value = 0                             # Z - PL_Ahh
    $ws.Cells.Item($rn, 27).Value = 0                             # AA - PL_Aha
}

Write-Output "Done."
